# Zensar Technologies Ltd - add "Exceptional items" column to the Quarterly
# Results sheet, between "P/l before exceptional items & tax" (col K) and
# "P/l before tax" (previously col L, now shifted to col M).
#
# This mirrors the authored change: a new blank column is inserted at L,
# pushing the existing L:T data/headers to M:U, and the new column gets the
# appropriate header text in the two header rows (row 1 = sentence case,
# row 2 = title case, matching the existing header-row casing convention
# used throughout this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new blank column at L; this shifts the existing columns L:T
# (and their values/styles) one position to the right, to M:U, and expands
# the used range from A1:T47 to A1:U47.
$ws.Range("L1:L47").EntireColumn.Insert()

# Populate the header cells for the newly-inserted column.
$ws.Range("L1").Value = "Exceptional items"
$ws.Range("L2").Value = "Exceptional Items"
